$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2:E2").NumberFormat = "@"
$ws.Range("G2").NumberFormat = "@"
$ws.Range("D2").Value = "308.15"
$ws.Range("E2").Value = "-1.07%"
$ws.Range("G2").Value = "8"

$ws.Range("D3:E3").NumberFormat = "@"
$ws.Range("G3").NumberFormat = "@"
$ws.Range("D3").Value = "37.58"
$ws.Range("E3").Value = "-0.19%"
$ws.Range("G3").Value = "8"

$ws.Range("D4:E4").NumberFormat = "@"
$ws.Range("G4").NumberFormat = "@"
$ws.Range("D4").Value = "5.137"
$ws.Range("E4").Value = "1.21%"
$ws.Range("G4").Value = "8"

$ws.Range("D5:E5").NumberFormat = "@"
$ws.Range("G5").NumberFormat = "@"
$ws.Range("D5").Value = "0.07838"
$ws.Range("E5").Value = "0.66%"
$ws.Range("G5").Value = "8"

$ws.Range("D6:E6").NumberFormat = "@"
$ws.Range("G6").NumberFormat = "@"
$ws.Range("D6").Value = "4.432"
$ws.Range("E6").Value = "1.75%"
$ws.Range("G6").Value = "8"

$ws.Range("B7:E7").NumberFormat = "@"
$ws.Range("G7").NumberFormat = "@"
$ws.Range("B7").Value = "KuCoinToken"
$ws.Range("C7").Value = "https://coinranking.com/coin/LOO6LmXd7G84Z+kucointoken-kcs"
$ws.Range("D7").Value = "8.280"
$ws.Range("E7").Value = "0.74%"
$ws.Range("G7").Value = "8"

$ws.Range("B8:E8").NumberFormat = "@"
$ws.Range("G8").NumberFormat = "@"
$ws.Range("B8").Value = "FTXToken"
$ws.Range("C8").Value = "https://coinranking.com/coin/NfeOYfNcl+ftxtoken-ftt"
$ws.Range("D8").Value = "1.901"
$ws.Range("E8").Value = "0.51%"
$ws.Range("G8").Value = "8"

$ws.Range("D9").NumberFormat = "@"
$ws.Range("G9").NumberFormat = "@"
$ws.Range("D9").Value = "2.992"
$ws.Range("G9").Value = "8"

$ws.Range("D10:E10").NumberFormat = "@"
$ws.Range("G10").NumberFormat = "@"
$ws.Range("D10").Value = "0.9284"
$ws.Range("E10").Value = "1.09%"
$ws.Range("G10").Value = "8"

$ws.Range("D11:E11").NumberFormat = "@"
$ws.Range("G11").NumberFormat = "@"
$ws.Range("D11").Value = "0.1084"
$ws.Range("E11").Value = "-9.36%"
$ws.Range("G11").Value = "8"

$ws.Range("D12:E12").NumberFormat = "@"
$ws.Range("G12").NumberFormat = "@"
$ws.Range("D12").Value = "0.1916"
$ws.Range("E12").Value = "-0.57%"
$ws.Range("G12").Value = "8"

$ws.Range("D13:E13").NumberFormat = "@"
$ws.Range("G13").NumberFormat = "@"
$ws.Range("D13").Value = "0.08966"
$ws.Range("E13").Value = "-1.53%"
$ws.Range("G13").Value = "8"

$ws.Range("D14:E14").NumberFormat = "@"
$ws.Range("G14").NumberFormat = "@"
$ws.Range("D14").Value = "0.03306"
$ws.Range("E14").Value = "-2.89%"
$ws.Range("G14").Value = "8"

$ws.Range("D15:E15").NumberFormat = "@"
$ws.Range("G15").NumberFormat = "@"
$ws.Range("D15").Value = "0.09587"
$ws.Range("E15").Value = "-1.22%"
$ws.Range("G15").Value = "8"

$ws.Range("D16:E16").NumberFormat = "@"
$ws.Range("G16").NumberFormat = "@"
$ws.Range("D16").Value = "0.001378"
$ws.Range("E16").Value = "1.12%"
$ws.Range("G16").Value = "8"

$ws.Range("B17:E17").NumberFormat = "@"
$ws.Range("G17").NumberFormat = "@"
$ws.Range("B17").Value = "CoinExToken"
$ws.Range("C17").Value = "https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet"
$ws.Range("D17").Value = "0.04386"
$ws.Range("E17").Value = "0.27%"
$ws.Range("G17").Value = "8"

$ws.Range("B18:E18").NumberFormat = "@"
$ws.Range("G18").NumberFormat = "@"
$ws.Range("B18").Value = "TigerCash"
$ws.Range("C18").Value = "https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("D18").Value = "0.005728"
$ws.Range("E18").Value = "-2.39%"
$ws.Range("G18").Value = "8"

$ws.Range("B19:E19").NumberFormat = "@"
$ws.Range("G19").NumberFormat = "@"
$ws.Range("B19").Value = "LEO"
$ws.Range("C19").Value = "https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("D19").Value = "3.614"
$ws.Range("E19").Value = "1.64%"
$ws.Range("G19").Value = "8"

$ws.Range("B20:E20").NumberFormat = "@"
$ws.Range("G20").NumberFormat = "@"
$ws.Range("B20").Value = "BitpandaEcosystemToken"
$ws.Range("C20").Value = "https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("D20").Value = "0.3455"
$ws.Range("E20").Value = "1.30%"
$ws.Range("G20").Value = "8"

$ws.Range("B21:E21").NumberFormat = "@"
$ws.Range("G21").NumberFormat = "@"
$ws.Range("B21").Value = "MCDex"
$ws.Range("C21").Value = "https://coinranking.com/coin/3nMM61qeg+mcdex-mcb"
$ws.Range("D21").Value = "6.320"
$ws.Range("E21").Value = "23.16%"
$ws.Range("G21").Value = "8"

$ws.Range("B22:E22").NumberFormat = "@"
$ws.Range("G22").NumberFormat = "@"
$ws.Range("B22").Value = "ProBitToken"
$ws.Range("C22").Value = "https://coinranking.com/coin/lQP4d6T2+probittoken-prob"
$ws.Range("D22").Value = "0.1276"
$ws.Range("E22").Value = "0.64%"
$ws.Range("G22").Value = "8"

$ws.Range("B23:E23").NumberFormat = "@"
$ws.Range("G23").NumberFormat = "@"
$ws.Range("B23").Value = "ZBToken"
$ws.Range("C23").Value = "https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb"
$ws.Range("D23").Value = "0.2583"
$ws.Range("E23").Value = "-0.08%"
$ws.Range("G23").Value = "8"

$ws.Range("D24:E24").NumberFormat = "@"
$ws.Range("G24").NumberFormat = "@"
$ws.Range("D24").Value = "0.001231"
$ws.Range("E24").Value = "1.63%"
$ws.Range("G24").Value = "8"

$ws.Range("D25:E25").NumberFormat = "@"
$ws.Range("G25").NumberFormat = "@"
$ws.Range("D25").Value = "0.004554"
$ws.Range("E25").Value = "6.67%"
$ws.Range("G25").Value = "8"

$ws.Range("D26:E26").NumberFormat = "@"
$ws.Range("G26").NumberFormat = "@"
$ws.Range("D26").Value = "0.0001301"
$ws.Range("E26").Value = "0.14%"
$ws.Range("G26").Value = "8"

$ws.Range("G27").NumberFormat = "@"
$ws.Range("G27").Value = "8"

$ws.Range("G28").NumberFormat = "@"
$ws.Range("G28").Value = "8"

$ws.Range("G29").NumberFormat = "@"
$ws.Range("G29").Value = "8"

$ws.Range("G30").NumberFormat = "@"
$ws.Range("G30").Value = "8"

$ws.Range("G31").NumberFormat = "@"
$ws.Range("G31").Value = "8"

$ws.Range("G32").NumberFormat = "@"
$ws.Range("G32").Value = "8"

$ws.Range("G33").NumberFormat = "@"
$ws.Range("G33").Value = "8"

$ws.Range("G34").NumberFormat = "@"
$ws.Range("G34").Value = "8"

$ws.Range("G35").NumberFormat = "@"
$ws.Range("G35").Value = "8"

$ws.Range("G36").NumberFormat = "@"
$ws.Range("G36").Value = "8"

$ws.Range("G37").NumberFormat = "@"
$ws.Range("G37").Value = "8"

$ws.Range("G38").NumberFormat = "@"
$ws.Range("G38").Value = "8"

$ws.Range("D39:E39").NumberFormat = "@"
$ws.Range("G39").NumberFormat = "@"
$ws.Range("D39").Value = "0.02202"
$ws.Range("E39").Value = "3.90%"
$ws.Range("G39").Value = "8"

$ws.Range("D40:E40").NumberFormat = "@"
$ws.Range("G40").NumberFormat = "@"
$ws.Range("D40").Value = "0.05017"
$ws.Range("E40").Value = "0.91%"
$ws.Range("G40").Value = "8"

$ws.Range("D41:E41").NumberFormat = "@"
$ws.Range("G41").NumberFormat = "@"
$ws.Range("D41").Value = "0.007444"
$ws.Range("E41").Value = "-4.55%"
$ws.Range("G41").Value = "8"

$ws.Range("D42:E42").NumberFormat = "@"
$ws.Range("G42").NumberFormat = "@"
$ws.Range("D42").Value = "0.1348"
$ws.Range("E42").Value = "-0.02%"
$ws.Range("G42").Value = "8"

$ws.Range("D43:E43").NumberFormat = "@"
$ws.Range("G43").NumberFormat = "@"
$ws.Range("D43").Value = "0.008719"
$ws.Range("E43").Value = "-11.99%"
$ws.Range("G43").Value = "8"

$ws.Range("D44:E44").NumberFormat = "@"
$ws.Range("G44").NumberFormat = "@"
$ws.Range("D44").Value = "0.002112"
$ws.Range("E44").Value = "2.56%"
$ws.Range("G44").Value = "8"

$ws.Range("D45:E45").NumberFormat = "@"
$ws.Range("G45").NumberFormat = "@"
$ws.Range("D45").Value = "0.007997"
$ws.Range("E45").Value = "-9.21%"
$ws.Range("G45").Value = "8"

$ws.Range("D46:E46").NumberFormat = "@"
$ws.Range("G46").NumberFormat = "@"
$ws.Range("D46").Value = "0.00006564"
$ws.Range("E46").Value = "-1.52%"
$ws.Range("G46").Value = "8"

$ws.Range("E47").NumberFormat = "@"
$ws.Range("G47").NumberFormat = "@"
$ws.Range("E47").Value = "0.13%"
$ws.Range("G47").Value = "8"

$ws.Range("E48").NumberFormat = "@"
$ws.Range("G48").NumberFormat = "@"
$ws.Range("E48").Value = "-1.83%"
$ws.Range("G48").Value = "8"

$ws.Range("D49:E49").NumberFormat = "@"
$ws.Range("G49").NumberFormat = "@"
$ws.Range("D49").Value = "0.001001"
$ws.Range("E49").Value = "-16.55%"
$ws.Range("G49").Value = "8"

$ws.Range("D50:E50").NumberFormat = "@"
$ws.Range("G50").NumberFormat = "@"
$ws.Range("D50").Value = "0.00002102"
$ws.Range("E50").Value = "0.13%"
$ws.Range("G50").Value = "8"

$ws.Range("D51:E51").NumberFormat = "@"
$ws.Range("G51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0002002"
$ws.Range("E51").Value = "0.13%"
$ws.Range("G51").Value = "8"
